$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.760.74"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.890.64"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'0.7914"
$ws.Range("E5").Value = "  -3.52%  "
$ws.Range("D6").Value = "'241.55"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D9").Value = "'25.54"
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("D10").Value = "'0.07034"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'0.08053"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "'0.7628"
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("D13").Value = "1.921.00"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "'5.308"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "29.771.67"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "'5.934"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "'243.21"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "'0.000007690"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.145.01"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'8.134"
$ws.Range("E23").Value = "  +16.70%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'0.1638"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").Value = "'9.306"
$ws.Range("D27").Value = "'163.68"
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("D28").Value = "'18.64"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").Value = "'2.049"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").Value = "'1.384"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "'1.534"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").Value = "'4.427"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").Value = "'0.05709"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("D34").Value = "'4.079"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "'1.263"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("D36").Value = "'0.7370"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'2.606"
$ws.Range("E38").Value = "  -3.78%  "
$ws.Range("D39").Value = "'0.01909"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").Value = "'2.776"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").Value = "'0.4399"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").Value = "'72.45"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "'5.831"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.8398"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "1.025.73"
$ws.Range("E46").Value = "  +3.87%  "
$ws.Range("D47").Value = "'102.28"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").Value = "'9.899"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").Value = "'1.845"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").Value = "'7.449"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").Value = "2.040.00"
$ws.Range("E51").Value = "  -0.98%  "